$wb = $excel.ActiveWorkbook

# --- Update SignUp 20 test data values (phone numbers / OTP codes) ---

# Login sheet: phone number used to kick off sign-up from the Login screen.
# This cell is formatted as "text entered with a leading apostrophe"
# (numFmtId 1 + quotePrefix), so re-assert the quote prefix to keep it
# stored as a shared string instead of being coerced to a number.
$wsLogin = $wb.Worksheets.Item("Login")
$wsLogin.Range("B12").Value = "'0363697350"

# SignUp sheet: phone number + OTP codes (cells are already Text-formatted).
$wsSignUp = $wb.Worksheets.Item("SignUp")
$wsSignUp.Range("B2").Value = "0363285179"
$wsSignUp.Range("E4").Value = "417116"
$wsSignUp.Range("E5").Value = "417116"
$wsSignUp.Range("E6").Value = "805154"
$wsSignUp.Range("E7").Value = "805154"
$wsSignUp.Range("E8").Value = "668334"

# --- "test app on browser stack": switch the active tab from
# ForgotPassword over to SignUp, landing on B8 ---
$wsSignUp.Activate()
$wsSignUp.Range("B8").Select()
